$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '64.517.48'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.29%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.325.55'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +0.62%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.998'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.42%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '528.33'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.66%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '175.93'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -3.64%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.589'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -2.57%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.330.90'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.93%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.999'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -0.24%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.607'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -1.76%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '53.98'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -10.11%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.135'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.79%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000259'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -0.70%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '8.99'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.60%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.832.13'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -0.31%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.118'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.30%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.303.07'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.22%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '64.336.93'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.50%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '17.43'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -1.62%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '11.21'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.25%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.958'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.50%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '382.30'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.61%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.16'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +6.96%  '
$ws.Range('B24').Value = 'RenderToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.27'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +1.40%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '81.41'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.24%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.71'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -3.22%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '6.11'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -0.90%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.73'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +1.03%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '11.29'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -2.80%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '8.23'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -2.60%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '28.97'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.15%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '632.47'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -2.78%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '6.68'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -1.01%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '11.21'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -1.16%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.106'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '57.51'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -2.25%  '
$ws.Range('E37').Value = '  +0.14%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '36.32'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.02%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.381'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -3.40%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0₃0757'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +7.05%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.996'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.29%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.24'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +12.97%  '
$ws.Range('B43').Value = 'Fetch.AI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.64'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +5.71%  '
$ws.Range('B44').Value = 'Kaspa'
$ws.Range('C44').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.126'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.94%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.922.67'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.06%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.18'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +5.46%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0401'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.44%  '
$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.67'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.65%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.63'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -3.05%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '138.23'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +1.58%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.125'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.42%  '
